$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Add([System.Type]::Missing, $ws)
$ws2.Name = "Sheet2"

$xvals = @(100,200,300,400,500,600,700,800,900,1000,1100,1200,1300,1400)
$yvals = @(49,116,183,251,319,387,456,525,592,662,730,798,867,937)
for ($i = 0; $i -lt $xvals.Length; $i++) {
    $row = 3 + $i
    $ws2.Cells.Item($row, 1).Value = $xvals[$i]
    $ws2.Cells.Item($row, 2).Value = $yvals[$i]
}

$co = $ws2.ChartObjects().Add(300, 100, 400, 300)
$chart = $co.Chart
$chart.ChartType = 74
$s1 = $chart.SeriesCollection(1)
$s1.XValues = "=Sheet2!`$A`$3:`$A`$16"
$s1.Values = "=Sheet2!`$B`$3:`$B`$16"
if ($chart.SeriesCollection().Count -gt 1) {
    $chart.SeriesCollection(2).Delete()
}
Write-Host "done"
